# Formed the consolidated report:
# Recompute the "Absent" column (H) for each attendance row based on the
# "Real" column (E): a student is marked absent (1) when there was no
# real attendance recorded that day (E = 0), otherwise present (0).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 3) { $lastRow = 21 }

for ($r = 3; $r -le $lastRow; $r++) {
    $realValue = $ws.Cells.Item($r, 5).Value2   # Column E - Real
    if ($realValue -eq $null) { $realValue = 0 }

    if ([double]$realValue -eq 0) {
        $ws.Cells.Item($r, 8).Value = 1   # Column H - Absent
    } else {
        $ws.Cells.Item($r, 8).Value = 0   # Column H - Absent
    }
}
